$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2239496666666667
$ws.Range("H2").Value = 0.671849
$ws.Range("I2").Value = 0.4264743968982249
$ws.Range("J2").Value = 0.4264743968982249
$ws.Range("M2").Value = 24.75542533333333
$ws.Range("N2").Value = 74.26627599999999
$ws.Range("O2").Value = 0.7762421087066456
$ws.Range("P2").Value = 0.7762421087066456
$ws.Range("Q2").Value = 5.543969251591555
$ws.Range("R2").Value = 49.895723264324
$ws.Range("S2").Value = 0.331047385157673
$ws.Range("T2").Value = 0.331047385157673

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2239496666666667
$ws.Range("H3").Value = 0.671849
$ws.Range("I3").Value = 0.4264743968982249
$ws.Range("J3").Value = 0.4264743968982249
$ws.Range("M3").Value = 3.818542
$ws.Range("N3").Value = 11.455626
$ws.Range("O3").Value = 0.1197358984688377
$ws.Range("P3").Value = 0.1197358984688377
$ws.Range("Q3").Value = 0.8551612080526668
$ws.Range("R3").Value = 7.696450872474001
$ws.Range("S3").Value = 0.05106429508656465
$ws.Range("T3").Value = 0.05106429508656464

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2239496666666667
$ws.Range("H4").Value = 0.671849
$ws.Range("I4").Value = 0.4264743968982249
$ws.Range("J4").Value = 0.4264743968982249
$ws.Range("M4").Value = 3.317404
$ws.Range("N4").Value = 9.952212
$ws.Range("O4").Value = 0.1040219928245168
$ws.Range("P4").Value = 0.1040219928245168
$ws.Range("Q4").Value = 0.7429315199986667
$ws.Range("R4").Value = 6.686383679988
$ws.Range("S4").Value = 0.04436271665398727
$ws.Range("T4").Value = 0.04436271665398727

# Row 5
$ws.Range("I5").Value = 0.4001470143891285
$ws.Range("J5").Value = 0.4001470143891285
$ws.Range("M5").Value = 24.75542533333333
$ws.Range("N5").Value = 74.26627599999999
$ws.Range("O5").Value = 0.7762421087066456
$ws.Range("P5").Value = 0.7762421087066456
$ws.Range("Q5").Value = 5.201725496358222
$ws.Range("R5").Value = 46.815529467224
$ws.Range("S5").Value = 0.3106109622420855
$ws.Range("T5").Value = 0.3106109622420855

# Row 6
$ws.Range("I6").Value = 0.4001470143891285
$ws.Range("J6").Value = 0.4001470143891285
$ws.Range("M6").Value = 3.818542
$ws.Range("N6").Value = 11.455626
$ws.Range("O6").Value = 0.1197358984688377
$ws.Range("P6").Value = 0.1197358984688377
$ws.Range("Q6").Value = 0.8023698649026667
$ws.Range("R6").Value = 7.221328784124
$ws.Range("S6").Value = 0.04791196228750522
$ws.Range("T6").Value = 0.04791196228750522

# Row 7
$ws.Range("I7").Value = 0.4001470143891285
$ws.Range("J7").Value = 0.4001470143891285
$ws.Range("M7").Value = 3.317404
$ws.Range("N7").Value = 9.952212
$ws.Range("O7").Value = 0.1040219928245168
$ws.Range("P7").Value = 0.1040219928245168
$ws.Range("Q7").Value = 0.6970684096986666
$ws.Range("R7").Value = 6.273615687287999
$ws.Range("S7").Value = 0.04162408985953773
$ws.Range("T7").Value = 0.04162408985953774

# Row 8
$ws.Range("G8").Value = 0.09104433333333334
$ws.Range("H8").Value = 0.273133
$ws.Range("I8").Value = 0.1733785887126465
$ws.Range("J8").Value = 0.1733785887126465
$ws.Range("M8").Value = 24.75542533333333
$ws.Range("N8").Value = 74.26627599999999
$ws.Range("O8").Value = 0.7762421087066456
$ws.Range("P8").Value = 0.7762421087066456
$ws.Range("Q8").Value = 2.253841195856444
$ws.Range("R8").Value = 20.284570762708
$ws.Range("S8").Value = 0.134583761306887
$ws.Range("T8").Value = 0.134583761306887

# Row 9
$ws.Range("G9").Value = 0.09104433333333334
$ws.Range("H9").Value = 0.273133
$ws.Range("I9").Value = 0.1733785887126465
$ws.Range("J9").Value = 0.1733785887126465
$ws.Range("M9").Value = 3.818542
$ws.Range("N9").Value = 11.455626
$ws.Range("O9").Value = 0.1197358984688377
$ws.Range("P9").Value = 0.1197358984688377
$ws.Range("Q9").Value = 0.3476566106953334
$ws.Range("R9").Value = 3.128909496258
$ws.Range("S9").Value = 0.02075964109476781
$ws.Range("T9").Value = 0.02075964109476781

# Row 10
$ws.Range("G10").Value = 0.09104433333333334
$ws.Range("H10").Value = 0.273133
$ws.Range("I10").Value = 0.1733785887126465
$ws.Range("J10").Value = 0.1733785887126465
$ws.Range("M10").Value = 3.317404
$ws.Range("N10").Value = 9.952212
$ws.Range("O10").Value = 0.1040219928245168
$ws.Range("P10").Value = 0.1040219928245168
$ws.Range("Q10").Value = 0.3020308355773333
$ws.Range("R10").Value = 2.718277520196
$ws.Range("S10").Value = 0.01803518631099176
$ws.Range("T10").Value = 0.01803518631099176

